# "Modify Fixed O&M for early retirement"
# Collapse the three year-keyed FOM-override rows (2018/2020/2025) on the
# INS sheet into a single always-on row: row 4 keeps its qualifiers (tech
# set + attribute) but the Year cell is dropped and every region flag
# (H:AH) becomes 1 ("apply"); rows 5 and 6 are wiped back to blank
# (formatted-but-empty) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# Row 4: drop the Year (E4) value, flip every regional FOM flag H4:AH4 from 0 to 1
$ws.Range("E4").ClearContents()
$ws.Range("H4:AH4").Value = 1

# Row 5: clear the whole "2020" scenario row back to empty (formatting stays)
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5:AH5").ClearContents()
$ws.Range("AJ5:AK5").ClearContents()

# Row 6: clear the whole "2025" scenario row back to empty (formatting stays)
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6:AH6").ClearContents()
$ws.Range("AJ6:AK6").ClearContents()

# Refresh the saved view/selection state on the INS sheet
$ws.Activate()
$appWindow = $excel.ActiveWindow
$appWindow.ScrollColumn = 20
$appWindow.ScrollRow = 1
$ws.Range("AG17").Select()
